# Applies the "total bcq nomination" update:
#  - Remove the PEDC column (old column E), shifting TOTAL_BCQ into column E
#  - Update SCPC (B), EDC (D) and TOTAL_BCQ (E) values for each hour row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the PEDC column (old column E). This shifts the old TOTAL_BCQ column (F)
# left into column E, matching the new dimension A1:E25 and header layout.
$ws.Range("E1").EntireColumn.Delete()

# New data values per hour (rows 2-25), columns B (SCPC), C (KSPC), D (EDC), E (TOTAL_BCQ)
$data = @(
    @(1,  25000, 20000, 12000, 57000),
    @(2,  25000, 20000, 0,     45000),
    @(3,  25000, 20000, 0,     45000),
    @(4,  25000, 20000, 0,     45000),
    @(5,  25000, 20000, 0,     45000),
    @(6,  25000, 20000, 0,     45000),
    @(7,  12500, 10000, 0,     22500),
    @(8,  25000, 20000, 0,     45000),
    @(9,  25000, 20000, 0,     45000),
    @(10, 25000, 20000, 20000, 65000),
    @(11, 25000, 20000, 20000, 65000),
    @(12, 25000, 20000, 20000, 65000),
    @(13, 25000, 20000, 20000, 65000),
    @(14, 25000, 20000, 20000, 65000),
    @(15, 25000, 20000, 20000, 65000),
    @(16, 25000, 20000, 20000, 65000),
    @(17, 25000, 20000, 20000, 65000),
    @(18, 25000, 20000, 20000, 65000),
    @(19, 25000, 20000, 20000, 65000),
    @(20, 25000, 20000, 20000, 65000),
    @(21, 25000, 20000, 20000, 65000),
    @(22, 25000, 20000, 20000, 65000),
    @(23, 25000, 20000, 20000, 65000),
    @(24, 25000, 20000, 20000, 65000)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
